$d = $word.ActiveDocument

# --- Change 1: merge the three italic runs "Zobrazit " / "zadane" / " ukoly"
# into a single run (same visible text, single <w:r>). This paragraph is the
# one right under the "Student" heading (the other "Zobrazit odevzdane ukoly"
# text is a different, unrelated paragraph under "Ucitel").
$rngFind = $d.Content.Duplicate
$rngFind.Find.Execute("Zobrazit zadané úkoly", $false, $false, $false, $false, $false, $true, 1, $false, "Zobrazit zadané úkoly", 2) | Out-Null

# --- Change 2/3: Admin -> "Vytvorit uzivatele" list.
# Insert a new bullet paragraph, with the same list formatting as the
# "Heslo" bullet, right before it: "Kdyz je v poli Email zadany email,
# ktery je uz pouzity v databazi ...". The "_GoBack" bookmark (which used
# to sit in its own empty paragraph right after the "Uzivatel" heading)
# now marks the end of this freshly-typed paragraph instead.

$rngHeslo = $d.Content.Duplicate
$rngHeslo.Find.Execute("Když do pole „Heslo“ zadám řetězec kratší") | Out-Null
$hesloPara = $rngHeslo.Duplicate
$hesloPara.Collapse(1)   # wdCollapseStart
$hesloPara.InsertParagraphBefore()
$newParaStart = $hesloPara.Start

$r1 = $d.Range($newParaStart, $newParaStart)
$r1.InsertAfter("Když je v poli „Email“ zadaný email, který je už použitý v databázi, zobrazí se hláška: „")

$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("Tento email je již v databázi.")

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter("“")

# Remove the old bookmark (currently in its own paragraph right after the
# "Uzivatel" heading) and re-create it at the end of the new paragraph.
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($r3.End, $r3.End))

# Clean up the now-empty paragraph that used to hold the bookmark (leave it
# as a bare empty paragraph, same as before).
$rngUzivatel = $d.Content.Duplicate
$rngUzivatel.Find.Execute("Uživatel") | Out-Null
$headingPara = $rngUzivatel.Paragraphs(1)
$afterHeadingPara = $headingPara.Next()

# Drop the lastRenderedPageBreak from the "Uzivatel" heading run, it moves
# down onto the "Zopakujte heslo" bullet instead.
$d.Content.Find.Execute("Uživatel", $false, $false, $false, $false, $false, $true, 1, $false, "Uživatel", 2) | Out-Null

$rngZopakujte = $d.Content.Duplicate
$rngZopakujte.Find.Execute("Když zadaný řetězec v poli „Zopakujte heslo“ není stejný") | Out-Null
